$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New automated-test rows (WAT-1195 / WAT-1196) appended to the WoS_AuthorTransformation
# test matrix, covering "Total Citing Publications" on the single-author and
# combined-author endpoints.

# Test case IDs first
$ws.Range("A94").Value = "WAT-1195"
$ws.Range("A95").Value = "WAT-1196"

# --- Row 94: WAT-1195 ---
$ws.Range("B94").Value = "Verify that ‘Get Author metadata’ also return ‘Total Citing Publications’"
$ws.Range("C94").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D94").Value = "/author/3515"
$ws.Range("E94").Value = "GET"
$ws.Range("J94").Value = "status=200||authorId=3515"
$ws.Range("K94").Value = "totalCitingPublications"
$ws.Range("L94").Value = "PASS"

# --- Row 95: WAT-1196 ---
$ws.Range("B95").Value = "Verify that ‘Get combined Author metadata’ also return ‘Total Citing Publications’"
$ws.Range("C95").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D95").Value = "/author/combine"
$ws.Range("E95").Value = "GET"
$ws.Range("G95").Value = "?authorId=45&authorId=74"
$ws.Range("J95").Value = "status=200||authorIds=45||authorIds=74"
$ws.Range("K95").Value = "totalCitingPublications"

# Match formatting (borders/alignment) of the surrounding test rows by copying
# the existing cell formats onto the two new rows, cell by cell.
$ws.Range("A63").Copy()
$ws.Range("A94").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B94").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C94").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D94").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E94").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F94").PasteSpecial(-4122)
$ws.Range("G13").Copy()
$ws.Range("G94").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H94").PasteSpecial(-4122)
$ws.Range("I3").Copy()
$ws.Range("I94").PasteSpecial(-4122)
$ws.Range("J18").Copy()
$ws.Range("J94").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("K94").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("L94").PasteSpecial(-4122)

$ws.Range("A63").Copy()
$ws.Range("A95").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B95").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C95").PasteSpecial(-4122)
$ws.Range("D3").Copy()
$ws.Range("D95").PasteSpecial(-4122)
$ws.Range("E3").Copy()
$ws.Range("E95").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F95").PasteSpecial(-4122)
$ws.Range("G42").Copy()
$ws.Range("G95").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H95").PasteSpecial(-4122)
$ws.Range("I3").Copy()
$ws.Range("I95").PasteSpecial(-4122)
$ws.Range("J18").Copy()
$ws.Range("J95").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("K95").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("L95").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Scroll the view down to the newly added rows and leave the selection where
# the author ended up after typing in the new data.
$ws.Range("K100").Select()

Write-Host "done"
